$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.415.24"
$ws.Range("E2").Value = "  +5.15%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.817.77"
$ws.Range("E3").Value = "  +5.86%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "343.96"
$ws.Range("E5").Value = "  +2.98%  "
$ws.Range("E6").Value = "  +0.17%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3813"
$ws.Range("E7").Value = "  +3.43%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3497"
$ws.Range("E8").Value = "  +4.45%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "48.77"
$ws.Range("E9").Value = "  -1.09%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.237"
$ws.Range("E10").Value = "  +4.33%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07740"
$ws.Range("E11").Value = "  +3.70%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.002"
$ws.Range("E12").Value = "  +0.07%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.11"
$ws.Range("E13").Value = "  +10.21%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.648"
$ws.Range("E14").Value = "  +5.61%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.246"
$ws.Range("E15").Value = "  +4.55%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.816.61"
$ws.Range("E16").Value = "  +5.58%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001122"
$ws.Range("E17").Value = "  +4.25%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06731"
$ws.Range("E18").Value = "  +1.53%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "86.16"
$ws.Range("E19").Value = "  +5.26%  "
$ws.Range("E20").Value = "  +0.15%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.63"
$ws.Range("E21").Value = "  +7.53%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.582"
$ws.Range("E22").Value = "  +8.14%  "
$ws.Range("E23").Value = "  +2.22%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "27.402.78"
$ws.Range("E24").Value = "  +5.30%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.466"
$ws.Range("E25").Value = "  -0.25%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.683"
$ws.Range("E26").Value = "  +9.38%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.03"
$ws.Range("E27").Value = "  +14.71%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.479"
$ws.Range("E28").Value = "  +12.10%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "154.37"
$ws.Range("E29").Value = "  +2.57%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.013.07"
$ws.Range("E30").Value = "  +5.28%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "135.92"
$ws.Range("E31").Value = "  +5.11%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.338"
$ws.Range("E32").Value = "  +7.09%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.026"
$ws.Range("E33").Value = "  -2.30%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "13.93"
$ws.Range("E34").Value = "  +8.17%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.08740"
$ws.Range("E35").Value = "  +2.41%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.708"
$ws.Range("E36").Value = "  -0.76%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.624"
$ws.Range("E37").Value = "  +5.19%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.7009"
$ws.Range("E38").Value = "  +13.68%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2281"
$ws.Range("E39").Value = "  +6.84%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.02422"
$ws.Range("E40").Value = "  +5.58%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.06489"
$ws.Range("E41").Value = "  +4.24%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.986"
$ws.Range("E42").Value = "  +5.12%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.298"
$ws.Range("E43").Value = "  +6.19%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.65"
$ws.Range("E44").Value = "  +1.77%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6554"
$ws.Range("E45").Value = "  +11.28%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.001"
$ws.Range("E46").Value = "  +0.03%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.052"
$ws.Range("E47").Value = "  +5.66%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.185"
$ws.Range("E48").Value = "  +8.22%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "133.10"
$ws.Range("E49").Value = "  +3.96%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07333"
$ws.Range("E50").Value = "  +0.90%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "80.72"
$ws.Range("E51").Value = "  +4.84%  "
